# Fix audit revise page error:
#  - header row (C1:F1) switches from the "1（%）".."4（%）" shared-string
#    labels to plain numbers 1..4 (keeps the existing header style)
#  - G2:G5 "overall" column switches from a bare percent number (60, 100,
#    0, 50) to a literal text percentage string ("60.0%", "100.0%", …)
#  - a handful of the C/F helper columns get recomputed numbers
#
# NOTE: Excel's normal "smart" cell-entry (Range.Value / .Formula / .Replace)
# auto-parses a string like "60.0%" back into the number 0.6 and stamps a
# percent NumberFormat style onto the cell - that's *not* what the source
# workbook has (it stores a literal shared-string "60.0%" with no style
# change). The Copy / PasteSpecial(xlPasteValues) round-trip below writes
# the text through a formula cell first, so the destination cell receives
# the literal string untouched by number/percent auto-detection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LiteralText($rangeAddress, $text) {
    # Stage the literal text in a scratch cell as a formula result (a
    # formula's cached string result is never re-parsed as a number/percent),
    # then copy just the value over to the real destination.
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy() | Out-Null
    $ws.Range($rangeAddress).PasteSpecial(-4163) | Out-Null
    $scratch.Clear() | Out-Null
}

# Row 1 header: C1:F1 become plain numbers (1,2,3,4) instead of the
# "1（%）".."4（%）" shared strings; style (s="1") stays as-is.
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
# G1 keeps showing "合格率" - value unchanged, just re-assert it.
$ws.Range("G1").Value = "合格率"

# Row 2 - 整体结果
$ws.Range("C2").Value = 0
$ws.Range("F2").Value = 3
Set-LiteralText "G2" "60.0%"

# Row 3 - 民族歧视
$ws.Range("F3").Value = 2
Set-LiteralText "G3" "100.0%"

# Row 4 - 信仰歧视
$ws.Range("C4").Value = 0
Set-LiteralText "G4" "0.0%"

# Row 5 - 国别歧视
$ws.Range("C5").Value = 0
$ws.Range("F5").Value = 1
Set-LiteralText "G5" "50.0%"
